$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 22, pushing the existing rows 22-26 down to 24-28.
$ws.Rows("22:23").Insert()

# New row 22: Comercializadora del Agro de Limarí - Damasco - Castle Brite - Especial
$ws.Range("A22").Value = 2
$ws.Range("B22").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 44902
$ws.Range("D22").Style = $ws.Range("D24").Style
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100103
$ws.Range("H22").Value = "Frutos de hueso (carozo)"
$ws.Range("I22").Value = 100103003
$ws.Range("J22").Value = "Damasco"
$ws.Range("K22").Value = "Castle Brite"
$ws.Range("L22").Value = "Especial"
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 25000
$ws.Range("O22").Value = 26000
$ws.Range("P22").Value = 25500
$ws.Range("Q22").Value = "$/caja 18 kilos"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 1417
$ws.Range("T22").Value = 18

# New row 23: Comercializadora del Agro de Limarí - Damasco - Castle Brite - Primera
$ws.Range("A23").Value = 2
$ws.Range("B23").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C23").Value = "Coquimbo"
$ws.Range("D23").Value = 44902
$ws.Range("D23").Style = $ws.Range("D24").Style
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100103
$ws.Range("H23").Value = "Frutos de hueso (carozo)"
$ws.Range("I23").Value = 100103003
$ws.Range("J23").Value = "Damasco"
$ws.Range("K23").Value = "Castle Brite"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 240
$ws.Range("N23").Value = 22000
$ws.Range("O23").Value = 23000
$ws.Range("P23").Value = 22500
$ws.Range("Q23").Value = "$/caja 18 kilos"
$ws.Range("R23").Value = "Región de O'Higgins"
$ws.Range("S23").Value = 1250
$ws.Range("T23").Value = 18
